$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.122.18'
$ws.Range("E2").Value = '  -5.89%  '

$ws.Range("D3").Value = '2.551.73'
$ws.Range("E3").Value = '  -1.85%  '

$ws.Range("E4").Value = '  -0.02%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '300.01'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -2.75%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '92.67'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -6.19%  '

$ws.Range("E7").Value = '  -3.19%  '

$ws.Range("E8").Value = '  +0.04%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.549'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -4.85%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '36.02'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  -6.70%  '

$ws.Range("E11").Value = '  -3.84%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '7.74'
$cell.Style = "Normal"

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.111'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +3.85%  '

$ws.Range("D14").Value = '2.948.09'
$ws.Range("E14").Value = '  -1.65%  '

$ws.Range("D15").Value = '2.564.33'
$ws.Range("E15").Value = '  -1.66%  '

$ws.Range("E16").Value = '  -4.04%  '

$ws.Range("E17").Value = '  -4.20%  '

$ws.Range("D18").Value = '43.147.09'
$ws.Range("E18").Value = '  -6.00%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '13.17'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +4.76%  '

$ws.Range("D20").Value = '0.0₃0981'
$ws.Range("E20").Value = '  -2.95%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '6.61'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -1.08%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '71.80'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -1.77%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '258.01'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -10.05%  '

$ws.Range("E24").Value = '  -2.76%  '

$ws.Range("B25").Value = 'ImmutableX'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.15'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -4.03%  '

$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '29.21'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +0.39%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -0.11%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '10.05'
$cell.Style = "Normal"

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '37.54'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -1.98%  '

$ws.Range("E30").Value = '  -3.18%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '5.98'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -4.71%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '154.16'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -2.35%  '

$ws.Range("B33").Value = 'WEMIXToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '2.76'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -1.51%  '

$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '2.17'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -3.33%  '

$ws.Range("E35").Value = '  -6.34%  '

$ws.Range("E36").Value = '  -4.03%  '

$ws.Range("E37").Value = '  -5.16%  '

$ws.Range("E38").Value = '  -2.06%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '17.06'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +8.96%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '23.37'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +10.19%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '3.45'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -2.34%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.0312'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -4.51%  '

$ws.Range("E43").Value = '  -2.29%  '

$ws.Range("D44").Value = '2.075.09'
$ws.Range("E44").Value = '  -1.30%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +0.00%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '85.06'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -9.73%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '8.95'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -2.85%  '

$ws.Range("E48").Value = '  +1.62%  '

$ws.Range("D49").Value = '2.803.02'
$ws.Range("E49").Value = '  -1.67%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '104.73'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -3.48%  '

$ws.Range("E51").Value = '  -3.17%  '
